$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8 reuses the same date-number-format as the cells above it (style index 1 / "m/d/yyyy"),
# so copy that formatting first, then fill in the values.
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("A8").Value = "Joshua Thomas Brooks"
$ws.Range("B8").Value = 45509
$ws.Range("C8").Value = 1632
$ws.Range("D8").Value = 33
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 130
$ws.Range("H8").Value = 79
$ws.Range("I8").Value = 105
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 97.7
$ws.Range("L8").Value = 96
$ws.Range("M8").Value = 127

$ws.Range("F9").Select()
